# Google Meet Attendance 04 jan 2024
# Adds the 04-Jan-2024 (H) attendance column and Admin's remark comments.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New date column header (H1 = 04-Jan-2024), matching the existing
# date header style used by the other day columns (G1 etc.) ---
$ws.Range("H1").NumberFormat = "d-mmm-yy"
$ws.Range("H1").Value = 45295

# --- Attendance entries for the new day ---
$ws.Range("H2").Value = "Reason"
$ws.Range("H3").Value = "Reason"
$ws.Range("H4").Value = "Present"
$ws.Range("H5").Value = "Present"
$ws.Range("H6").Value = "Present"
$ws.Range("H7").Value = "Reason"
$ws.Range("H8").Value = "Present"
$ws.Range("H9").Value = "Absent"

# --- Admin's remark comments for the "Reason"/"Absent" entries ---
function Add-AdminComment {
    param(
        [string]$CellRef,
        [string]$FullText
    )

    $cell = $ws.Range($CellRef)
    $cmt = $cell.AddComment($FullText)
    return $cmt
}

Add-AdminComment "H2" "Admin:`ntravelling in train `nReached late from office"
Add-AdminComment "H3" "Admin:`nUnavailable for Meet"
Add-AdminComment "H7" "Admin:`nRelocating the Home"
Add-AdminComment "H9" "Admin:`nTravelling"

# --- Selection state left by the editor (D1:G2 highlighted, G2 active) ---
$null = $ws.Range("D1:G2").Select()
